$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")
$fcs = $ws.Range("H5:BJ40").FormatConditions
Write-Host $fcs.Count
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    Write-Host $fc.Type
}
